$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find.Execute failed for: $find"
    }
}

# 1) Ativação date bump
Replace-Text "Ativação: 01/01/2012" "Ativação: 01/01/2024"

# 2) Curso (semestre ideal) - drop EQD (6)
Replace-Text "Curso (semestre ideal): EA (6), EB (5), EQD (6), EQN (6)" "Curso (semestre ideal): EA (6), EB (5), EQN (6)"

# 3) Programa paragraph - merge the 7 numbered items into a single run (drop manual line breaks)
$programaFind = "1) Introdução à Análise Instrumental. Correlação entre métodos analíticos instrumentais e por via úmida. Preparo de amostras em meio sólido e em meios líquidos aquosos e não aquosos. Solubilização, digestão, fontes de energia aplicadas ao preparo, estabilização de amostras.^l2) Introdução aos Métodos Espectroanalíticos: Natureza da energia radiante. Espectro eletromagnético. Interação da radiação com a matéria. Absorção seletiva. Absortividade. Lei de Beer-Lambert. Curvas analíticas. ^l3) Introdução à Espectrofotometria no UV/Visível. Instrumentação. Aplicações e interpretação de resultados. Determinações simultâneas. Parte Experimental.^l4) Introdução às Espectrometrias de Absorção e de Emissão Atômicas. Instrumentação. Interferências. Origem do espectro de emissão atômica. Fontes de atomização e de excitação. Calibração. Aplicações e interpretação de resultados. Parte Experimental.^l5) Introdução à Espectroscopia no Infravermelho. Instrumentação. Interpretação de espectros. Aplicações. Parte Experimental.^l6) Introdução aos Métodos Eletroanalíticos: Potenciometria e Condutimetria. Instrumentação.  Métodos diretos e indiretos. Aplicações e interpretação de resultados. Parte experimental.^l7) Introdução aos Métodos Cromatográficos. Conceitos básicos dos métodos de separação. Fases móvel e estacionária. Cromatografia planar em papel e em camada delgada. Cromatografia em coluna: cromatografia a gás e cromatografia líquida de alta eficiência. Instrumentação. Aplicações e interpretação de resultados. Parte Experimental."
$programaReplace = "1) Introdução à Análise Instrumental. Correlação entre métodos analíticos instrumentais e por via úmida. Preparo de amostras em meio sólido e em meios líquidos aquosos e não aquosos. Solubilização, digestão, fontes de energia aplicadas ao preparo, estabilização de amostras.2) Introdução aos Métodos Espectroanalíticos: Natureza da energia radiante. Espectro eletromagnético. Interação da radiação com a matéria. Absorção seletiva. Absortividade. Lei de Beer-Lambert. Curvas analíticas. 3) Introdução à Espectrofotometria no UV/Visível. Instrumentação. Aplicações e interpretação de resultados. Determinações simultâneas. Parte Experimental.4) Introdução às Espectrometrias de Absorção e de Emissão Atômicas. Instrumentação. Interferências. Origem do espectro de emissão atômica. Fontes de atomização e de excitação. Calibração. Aplicações e interpretação de resultados. Parte Experimental.5) Introdução à Espectroscopia no Infravermelho. Instrumentação. Interpretação de espectros. Aplicações. Parte Experimental.6) Introdução aos Métodos Eletroanalíticos: Potenciometria e Condutimetria. Instrumentação.  Métodos diretos e indiretos. Aplicações e interpretação de resultados. Parte experimental.7) Introdução aos Métodos Cromatográficos. Conceitos básicos dos métodos de separação. Fases móvel e estacionária. Cromatografia planar em papel e em camada delgada. Cromatografia em coluna: cromatografia a gás e cromatografia líquida de alta eficiência. Instrumentação. Aplicações e interpretação de resultados. Parte Experimental."
Replace-Text $programaFind $programaReplace

# 4) Avaliação - Método text
Replace-Text "A avaliação da disciplina será feita por meio de avaliações escritas individuais (provas) e avaliações de atividades em grupo (relatórios das aulas práticas e/ou trabalhos escritos e/ou apresentações de seminários)." "A avaliação será feita por meio de duas provas (P1 e P2). A critério do professor, a avaliação poderá ser complementada por meio de trabalhos e/ou relatórios, valendo até 30% da nota das provas."

# 5) Avaliação - Critério text
Replace-Text "A Média Final (MF) será calculada pela média entre todas as avaliações realizadas durante o semestre, sendo o conjunto das avaliações individuais correspondentes a 75% da composição de MF e o conjunto das avaliações em grupo correspondentes a 25% da composição de MF. Será aprovado o aluno que obtiver MF maior ou igual a cinco e frequência mínima de 70% no semestre." "A nota final (NF) será calculada pela média aritmética das provas. NF=(P1 +P2)/2."

# 6) Avaliação - Norma de recuperação text (merges two runs/breaks into one new run)
$normaFind = "No período de Recuperação haverá horário previamente definido para resolução de dúvidas e será realizada uma avaliação escrita individual (Prova da Recuperação = PR), com conteúdo de todos os tópicos apresentados na disciplina durante o semestre.^lA Nota de Recuperação (NR) será dada pela média aritmética entre a Média do Semestre (MF) e a Prova da Recuperação (PR), sendo considerado aprovado o aluno que obtiver NR maior ou igual a cinco."
$normaReplace = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."
Replace-Text $normaFind $normaReplace

# 7) Bibliografia paragraph - drop/merge several references into a single run
$biblioFind = "1) Skoog, D.A.; Holler, F.J. ; Nieman, T.A. Princípios de análise instrumental. 5. ed. Porto Alegre: Bookman,  2002.^l2) MENDHAM,J.; DENNEY, R.C.; BARNES, J.D. ; Thomas, M. Vogel: análise química quantitativa. 6. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2002.^l3) OHLWEILER, O.A. Fundamentos de análise instrumental. Rio de Janeiro: Livros Técnicos e Científicos, 1981.^l4) KRUG, F.J. (org.) Métodos de preparo de amostras: fundamentos sobre métodos de preparo de amostras orgânicas e inorgânicas para análise elementar. 1. ed. Piracicaba: Edição do autor, 2008. ^l5) COLLINS, C.H.; BRAGA, G.L.; BONATO, P.S. (Org.) Fundamentos de cromatografia. 1. ed. Campinas: Editora da UNICAMP, 2006.^l^lBibliografia complementar^l1) CHRISTIAN, G.D. Analytical chemistry. 4. ed. Nova York: John Wiley & Sons, 1986.^l2) DYER, J.R. Aplicação da espectroscopia de absorção aos compostos orgânicos. 1. Reimpressão. São Paulo: Edgard Blucher, 1977.^l3) SILVERSTEIN, R.M.; WEBSTER, F.X.; KIEMLE, D.J. Identificação espectrométrica de compostos orgânicos. 7. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2007.^l4) WILLARD, H.H.; MERRITE, L.; DEAB, J. Instrumentação analítica. Lisboa: Fundação Calouste Gulbekian,  1989."
$biblioReplace = "1) Skoog, D.A.; Holler, F.J. ; Nieman, T.A. Princípios de análise instrumental. 6a. ed. Porto Alegre: Bookman, 2009.2) KRUG, F.J. (org.) Métodos de preparo de amostras: fundamentos sobre métodos de preparo de amostras orgânicas e inorgânicas para análise elementar. 1. ed. Piracicaba: Edição do autor, 2008.3) COLLINS, C.H.; BRAGA, G.L.; BONATO, P.S. (Org.) Fundamentos de cromatografia. 1. ed. Campinas: Editora da UNICAMP, 2006.Bibliografia complementar1) CHRISTIAN, G.D. Analytical chemistry. 4. ed. Nova York: John Wiley & Sons, 1986.3) SILVERSTEIN, R.M.; WEBSTER, F.X.; KIEMLE, D.J. Identificação espectrométrica de compostos orgânicos. 7. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2007.4) WILLARD, H.H.; MERRITE, L.; DEAB, J. Instrumentação analítica. Lisboa: Fundação Calouste Gulbekian, 1989."
Replace-Text $biblioFind $biblioReplace

Write-Output "All replacements applied."
